$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 6837
$ws.Range("B2").Value = "Rebeca Alves"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45096
$ws.Range("G2").Value = 11982.72

# Row 3
$ws.Range("A3").Value = 91004
$ws.Range("B3").Value = "Isabella Pires"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 3504.31

# Row 4
$ws.Range("A4").Value = 5782
$ws.Range("B4").Value = "Ryan da Costa"
$ws.Range("C4").Value = "P&D"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45099
$ws.Range("G4").Value = 3711.57

# Row 5
$ws.Range("A5").Value = 54491
$ws.Range("B5").Value = "Gabriel Ribeiro"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45104
$ws.Range("G5").Value = 3738.06

# Row 6
$ws.Range("A6").Value = 35051
$ws.Range("B6").Value = "Thiago Silva"
$ws.Range("C6").Value = "Marketing"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45082
$ws.Range("G6").Value = 8530.440000000001

# Row 7
$ws.Range("A7").Value = 66247
$ws.Range("B7").Value = "Lucca Campos"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45080
$ws.Range("G7").Value = 8603.280000000001

# Row 8
$ws.Range("A8").Value = 73519
$ws.Range("B8").Value = "Julia Gomes"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("F8").Value = 45078
$ws.Range("G8").Value = 4688.54

# Row 9
$ws.Range("A9").Value = 6829
$ws.Range("B9").Value = "Raul da Cruz"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45099
$ws.Range("G9").Value = 5294.9

# Row 10
$ws.Range("A10").Value = 95444
$ws.Range("B10").Value = "Srta. Ana Júlia Barros"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45091
$ws.Range("G10").Value = 6349.58

# Row 11
$ws.Range("A11").Value = 80367
$ws.Range("B11").Value = "Augusto Porto"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45091
$ws.Range("G11").Value = 5182.44
